$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (Mac-Address / Document type entries) at row 33,
# mirroring the existing table's structure.
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 10032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"

# Mirror the scroll/selection state saved with the edit.
$ws.Range("E29").Select()
